{"js": "// The paragraph \"Petra (\" + \"neutral disappointed\" + \"): Oh, I see\u2026\" is\n// currently split across three separate runs. The edit merges them into a\n// single run containing the full text \"Petra (neutral disappointed): Oh, I see\u2026\"\n// (formatting/content is otherwise unchanged).\nconst body = context.document.body;\nconst fullText = \"Petra (neutral disappointed): Oh, I see\\u2026\";\n\nconst results = body.search(fullText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target paragraph text: \" + fullText);\n}\n\n// Replacing the matched range's text with the same text collapses the\n// underlying multiple runs into a single run while preserving the\n// formatting of the range's leading run.\nresults.items[0].insertText(fullText, \"Replace\");\nawait context.sync();\n", "ps1": "# The paragraph \"Petra (\" + \"neutral disappointed\" + \"): Oh, I see\u2026\" is\n# currently split across three separate runs. This rewrites that span so\n# the text ends up in a single run reading\n# \"Petra (neutral disappointed): Oh, I see\u2026\" (content is unchanged, only\n# the run layout collapses to one run).\n$d = $word.ActiveDocument\n\n$ellipsis = [char]0x2026\n$target = \"Petra (neutral disappointed): Oh, I see\" + $ellipsis\n\n$rng = $d.Content\n$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $target, 2)\n\nif (-not $found) {\n    throw \"Could not find target paragraph text: $target\"\n}\n"}
